$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row 24 that is a copy of the current row 23
# (phone stored as text "09876543", blank birthday, 0 points).
$ws.Range("A24").Value = "'09876543"
$ws.Range("A24").Style = "Normal"
$ws.Range("B24").Value = "'"
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = 0

# Row 23's phone number becomes numeric (loses the leading zero): 09876543 -> 9876543
$ws.Range("A23").Value = 9876543
